$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-28 03:48:12"
$ws.Range("H2").Value = "'86%"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("N2").Value = "0.4 °C 3:04 TU"
$ws.Range("E3").Value = "2026-02-28 03:48:15"
$ws.Range("H3").Value = "'82%"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("O3").Value = "-0.2 °C"
$ws.Range("E4").Value = "2026-02-28 03:48:17"
$ws.Range("H4").Value = "'97%"
$ws.Range("C4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("L4").Value = "8.3 km/h - 334º 3:09 TU"
$ws.Range("O4").Value = "7.5 °C"
$ws.Range("E5").Value = "2026-02-28 03:48:19"
$ws.Range("N5").Value = "-0.3 °C 3:17 TU"
$ws.Range("E6").Value = "2026-02-28 03:48:22"
$ws.Range("E7").Value = "2026-02-28 03:48:24"
$ws.Range("J7").Value = "1023.7 hPa"
$ws.Range("E8").Value = "2026-02-28 03:48:26"
$ws.Range("H8").Value = "'96%"
$ws.Range("C8").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").Value = "1023.8 hPa"
$ws.Range("L8").Value = "14.8 km/h - 57º 3:04 TU"
$ws.Range("E9").Value = "2026-02-28 03:48:29"
$ws.Range("M9").Value = "8.4 °C 3:24 TU"
$ws.Range("O9").Value = "7.3 °C"
$ws.Range("E10").Value = "2026-02-28 03:48:31"
$ws.Range("O10").Value = "7.4 °C"
$ws.Range("E11").Value = "2026-02-28 03:48:33"
$ws.Range("H11").Value = "'93%"
$ws.Range("C11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("N11").Value = "3.1 °C 3:20 TU"
$ws.Range("O11").Value = "3.9 °C"
$ws.Range("E12").Value = "2026-02-28 03:48:36"
$ws.Range("E13").Value = "2026-02-28 03:48:38"
$ws.Range("H13").Value = "'84%"
$ws.Range("C13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("J13").Value = "1025.8 hPa"
$ws.Range("N13").Value = "-0.4 °C 3:09 TU"
$ws.Range("O13").Value = "1.9 °C"
$ws.Range("E14").Value = "2026-02-28 03:48:40"
$ws.Range("H14").Value = "'99%"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("M14").Value = "10.7 °C 3:29 TU"
$ws.Range("O14").Value = "9.9 °C"
$ws.Range("E15").Value = "2026-02-28 03:48:43"
$ws.Range("M15").Value = "7.9 °C 3:17 TU"
$ws.Range("O15").Value = "6.5 °C"
$ws.Range("E16").Value = "2026-02-28 03:48:45"
$ws.Range("H16").Value = "'64%"
$ws.Range("C16").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("N16").Value = "-1.4 °C 3:18 TU"
$ws.Range("O16").Value = "-0.5 °C"
$ws.Range("E17").Value = "2026-02-28 03:48:47"
$ws.Range("H17").Value = "'44%"
$ws.Range("C17").Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$ws.Range("M17").Value = "5.0 °C 3:26 TU"
$ws.Range("E18").Value = "2026-02-28 03:48:50"
$ws.Range("M18").Value = "8.6 °C 3:29 TU"
$ws.Range("O18").Value = "7.8 °C"
$ws.Range("E19").Value = "2026-02-28 03:48:52"
$ws.Range("O19").Value = "7.9 °C"
$ws.Range("E20").Value = "2026-02-28 03:48:54"
$ws.Range("L20").Value = "22.0 km/h - 176º 3:15 TU"
$ws.Range("O20").Value = "0.3 °C"
$ws.Range("E21").Value = "2026-02-28 03:48:57"
$ws.Range("N21").Value = "4.8 °C 3:26 TU"
$ws.Range("O21").Value = "5.7 °C"
$ws.Range("E22").Value = "2026-02-28 03:48:59"
$ws.Range("H22").Value = "'61%"
$ws.Range("C22").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("N22").Value = "-1.5 °C 3:22 TU"
$ws.Range("O22").Value = "-0.6 °C"
$ws.Range("E23").Value = "2026-02-28 03:49:01"
$ws.Range("H23").Value = "'68%"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null
$ws.Range("N23").Value = "-0.7 °C 3:23 TU"
$ws.Range("O23").Value = "0.3 °C"
$ws.Range("E24").Value = "2026-02-28 03:49:04"
$ws.Range("N24").Value = "4.3 °C 3:18 TU"
$ws.Range("O24").Value = "6.5 °C"
$ws.Range("E25").Value = "2026-02-28 03:49:06"
$ws.Range("O25").Value = "1.4 °C"
$ws.Range("E26").Value = "2026-02-28 03:49:08"
$ws.Range("H26").Value = "'72%"
$ws.Range("C26").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
$ws.Range("J26").Value = "1023.9 hPa"
$ws.Range("E27").Value = "2026-02-28 03:49:11"
$ws.Range("N27").Value = "1.6 °C 3:06 TU"
$ws.Range("O27").Value = "2.9 °C"
$ws.Range("E28").Value = "2026-02-28 03:49:13"
$ws.Range("N28").Value = "5.5 °C 3:04 TU"
$ws.Range("O28").Value = "6.6 °C"
$ws.Range("E29").Value = "2026-02-28 03:49:16"
$ws.Range("M29").Value = "9.7 °C 3:29 TU"
$ws.Range("O29").Value = "8.4 °C"
$ws.Range("E30").Value = "2026-02-28 03:49:18"
$ws.Range("E31").Value = "2026-02-28 03:49:20"
$ws.Range("J31").Value = "1023.7 hPa"
$ws.Range("L31").Value = "39.6 km/h - 2º 3:11 TU"
$ws.Range("N31").Value = "9.8 °C 3:13 TU"
$ws.Range("E32").Value = "2026-02-28 03:49:23"
$ws.Range("H32").Value = "'85%"
$ws.Range("C32").Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "2026-02-28 03:49:25"
$ws.Range("H33").Value = "'68%"
$ws.Range("C33").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null
$ws.Range("L33").Value = "9.0 km/h - 86º 3:28 TU"
$ws.Range("N33").Value = "4.1 °C 3:09 TU"
$ws.Range("O33").Value = "5.3 °C"
$ws.Range("E34").Value = "2026-02-28 03:49:27"
$ws.Range("N34").Value = "-0.7 °C 3:16 TU"
$ws.Range("O34").Value = "0.1 °C"
$ws.Range("E35").Value = "2026-02-28 03:49:30"
$ws.Range("L35").Value = "27.0 km/h - 258º 3:03 TU"
$ws.Range("N35").Value = "6.6 °C 3:29 TU"
$ws.Range("O35").Value = "7.1 °C"
$ws.Range("E36").Value = "2026-02-28 03:49:32"
$ws.Range("J36").Value = "1024.3 hPa"
$ws.Range("N36").Value = "9.4 °C 3:03 TU"
$ws.Range("E37").Value = "2026-02-28 03:49:34"
$ws.Range("H37").Value = "'86%"
$ws.Range("C37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("J37").Value = "1025.9 hPa"
$ws.Range("O37").Value = "4.8 °C"
$ws.Range("E38").Value = "2026-02-28 03:49:36"
$ws.Range("E39").Value = "2026-02-28 03:49:39"
$ws.Range("N39").Value = "-0.8 °C 3:26 TU"
$ws.Range("O39").Value = "0.6 °C"
$ws.Range("E40").Value = "2026-02-28 03:49:41"
$ws.Range("N40").Value = "3.3 °C 3:23 TU"
$ws.Range("E41").Value = "2026-02-28 03:49:43"
$ws.Range("J41").Value = "1023.0 hPa"
$ws.Range("E42").Value = "2026-02-28 03:49:46"
$ws.Range("O42").Value = "7.6 °C"
$ws.Range("E43").Value = "2026-02-28 03:49:48"
$ws.Range("E44").Value = "2026-02-28 03:49:50"
$ws.Range("H44").Value = "'90%"
$ws.Range("C44").Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4122) | Out-Null
$ws.Range("O44").Value = "-0.9 °C"
$ws.Range("E45").Value = "2026-02-28 03:49:52"
$ws.Range("H45").Value = "'87%"
$ws.Range("C45").Copy() | Out-Null
$ws.Range("H45").PasteSpecial(-4122) | Out-Null
$ws.Range("J45").Value = "1023.7 hPa"
$ws.Range("N45").Value = "6.6 °C 3:28 TU"
$ws.Range("O45").Value = "7.4 °C"
$ws.Range("E46").Value = "2026-02-28 03:49:55"
$ws.Range("N46").Value = "10.6 °C 3:17 TU"
$excel.CutCopyMode = 0
